$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a native Excel table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
